# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the latest generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")

$exhibitUpdates = @{
    2  = 7282
    5  = 26
    6  = 570
    7  = 195
    8  = 136
    11 = 63
    12 = 224
    14 = 468
    15 = 32
    16 = 1871
    18 = 48
    19 = 3808
    21 = 254
    23 = 45
    24 = 3
    25 = 38
    26 = 2490
    27 = 25
    28 = 319
    31 = 44
    33 = 25
    38 = 27
    39 = 1484
    40 = 163
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")

$allUpdates = @{
    2  = 7282
    5  = 26
    7  = 570
    8  = 195
    9  = 136
    12 = 63
    13 = 224
    15 = 468
    16 = 32
    17 = 1872
    19 = 48
    20 = 3808
    22 = 254
    24 = 45
    25 = 3
    26 = 38
    27 = 2490
    28 = 25
    29 = 319
    32 = 44
    34 = 25
    38 = 166
    39 = 27
    40 = 1484
    41 = 163
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
